$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = -0.08318890944582957
$ws.Range("C4").Value = -0.05790617081637017
$ws.Range("D4").Value = -0.6334156050371067
$ws.Range("E4").Value = -0.7451615853267117
$ws.Range("F4").Value = -0.1996804606122911
$ws.Range("G4").Value = -0.07295128576116339
$ws.Range("B7").Value = 0.6829905449383783
$ws.Range("C7").Value = 0.6761456942398975
$ws.Range("D7").Value = 0.7395152237330389
$ws.Range("E7").Value = 0.6926365792964611
$ws.Range("F7").Value = 0.6926107193993014
$ws.Range("G7").Value = 0.6852021530792638
$ws.Range("B8").Value = 0.680911794413384
$ws.Range("C8").Value = 0.6730383183012454
$ws.Range("D8").Value = 0.6722933459867264
$ws.Range("E8").Value = 0.6033466470250419
$ws.Range("F8").Value = 0.6813383339329853
$ws.Range("G8").Value = 0.6525936111671558
$ws.Range("B9").Value = 1.92704027726453
$ws.Range("C9").Value = 1.921723746690912
$ws.Range("D9").Value = 1.654475624169477
$ws.Range("E9").Value = 1.734853893503501
$ws.Range("F9").Value = 1.97305120111073
$ws.Range("G9").Value = 2.044630187052906
$ws.Range("C10").Value = 113
$ws.Range("D10").Value = 4
$ws.Range("F10").Value = 25
$ws.Range("G10").Value = 15
$ws.Range("B11").Value = 0.7591228491109451
$ws.Range("C11").Value = 0.7239977444526943
$ws.Range("D11").Value = 1.037871217480483
$ws.Range("E11").Value = 1.04494793762322
$ws.Range("F11").Value = 0.8498535038538788
$ws.Range("G11").Value = 0.6574432098520148
$ws.Range("B12").Value = 0.7641007038592136
$ws.Range("C12").Value = 0.7309444891176156
$ws.Range("D12").Value = 1.305708951023833
$ws.Range("E12").Value = 1.348508232351754
$ws.Range("F12").Value = 0.8810187945452763
$ws.Range("G12").Value = 0.7255448969283191
$ws.Range("B13").Value = -1.623443384313227
$ws.Range("C13").Value = -1.535011378627223
$ws.Range("D13").Value = -2.039262426559294
$ws.Range("E13").Value = -1.904088349504508
$ws.Range("F13").Value = -1.855251210397697
$ws.Range("G13").Value = -1.548606873020195
$ws.Range("C14").Value = -59
$ws.Range("D14").Value = -3
$ws.Range("F14").Value = -20
$ws.Range("G14").Value = -15
